$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cases")
$ws.Range("AA4").Value = 2
$ws.Range("AB4").Value = 5
$ws.Range("AA5").ClearContents()
$ws.Range("AA6").Value = 6
$ws.Range("AB6").Value = 18
$ws.Range("AA7").Value = 7
$ws.Range("AB7").Value = 25
$ws.Range("AA8").Value = 11
$ws.Range("AB8").Value = 35
$ws.Range("AA9").Value = 14
$ws.Range("AB9").Value = 51
$ws.Range("AA10").Value = 16
$ws.Range("AB10").Value = 70
$ws.Range("AA11").Value = 24
$ws.Range("AB11").Value = 114
$ws.Range("AA12").Value = 30
$ws.Range("AB12").Value = 180
$ws.Range("AA13").Value = 35
$ws.Range("AB13").Value = 252
$ws.Range("AA14").Value = 41
$ws.Range("AB14").Value = 307
$ws.Range("AA15").Value = 50
$ws.Range("AB15").Value = 383
$ws.Range("AA16").Value = 63
$ws.Range("AB16").Value = 472
$ws.Range("AA17").Value = 102
$ws.Range("AB17").Value = 594
$ws.Range("AA18").Value = 141
$ws.Range("AB18").Value = 901
$ws.Range("AA19").Value = 164
$ws.Range("AB19").Value = 1215
$ws.Range("AA20").Value = 219
$ws.Range("AB20").Value = 1595
$ws.Range("AA21").Value = 251
$ws.Range("AB21").Value = 1853
$ws.Range("AA22").Value = 327
$ws.Range("AB22").Value = 2407
$ws.Range("AA23").Value = 430
$ws.Range("AB23").Value = 3032
$ws.Range("AA24").Value = 569
$ws.Range("AB24").Value = 4110
$ws.Range("AA25").Value = 680
$ws.Range("AB25").Value = 5383
$ws.Range("AA26").Value = 712
$ws.Range("AB26").Value = 6508
$ws.Range("AB27").Value = 7321
$ws.Range("AB28").Value = 7929
$ws.Range("AA29").Value = 1076
$ws.Range("AB29").Value = 9175
$ws.Range("AA30").Value = 1224
$ws.Range("AB30").Value = 10169
$ws.Range("AA31").Value = 1371
$ws.Range("AB31").Value = 11196
$ws.Range("AA32").Value = 1503
$ws.Range("AB32").Value = 12451
$ws.Range("AA33").Value = 1630
$ws.Range("AB33").Value = 13747
$ws.Range("AA34").Value = 1704
$ws.Range("AB34").Value = 14661
$ws.Range("AA35").Value = 1736
$ws.Range("AB35").Value = 15407
$ws.Range("AA36").Value = 1862
$ws.Range("AB36").Value = 16443
$ws.Range("AA37").Value = 1954
$ws.Range("AB37").Value = 17429
$ws.Range("AA38").Value = 2143
$ws.Range("AB38").Value = 18450
$ws.Range("AA39").Value = 2307
$ws.Range("AB39").Value = 19548
$ws.Range("AA40").Value = 2435
$ws.Range("AB40").Value = 20465
$ws.Range("AA41").Value = 2468
$ws.Range("AB41").Value = 21064
$ws.Range("AA42").Value = 2498
$ws.Range("AB42").Value = 21563
$ws.Range("AA43").Value = 2612
$ws.Range("AB43").Value = 22187
$ws.Range("AA44").Value = 2696
$ws.Range("AB44").Value = 22818
$ws.Range("AA45").Value = 2791
$ws.Range("AB45").Value = 23603
$ws.Range("AA46").Value = 2888
$ws.Range("AB46").Value = 24249
$ws.Range("AA47").Value = 2928
$ws.Range("AB47").Value = 24744
$ws.Range("AA48").Value = 2986
$ws.Range("AB48").Value = 25229
$ws.Range("B49").Value = 899
$ws.Range("AA49").Value = 3003
$ws.Range("AB49").Value = 25555
$ws.Range("B50").Value = 906
$ws.Range("AA50").Value = 3020
$ws.Range("AB50").Value = 25753
$ws.Range("B51").Value = 912
$ws.Range("AB51").Value = 25931
$ws = $wb.Worksheets.Item("Fatalities")
$ws.Range("B49").Value = 18
$ws.Range("B50").Value = 19
$ws.Range("AB50").Value = 1143
$ws.Range("B51").Value = 19
$ws.Range("AB51").Value = 1161
$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Range("B51").Value = 78
$ws.Range("AB51").Value = 1900
$ws = $wb.Worksheets.Item("ICU")
$ws.Range("B51").Value = 22
$ws.Range("AB51").Value = 356
$ws = $wb.Worksheets.Item("Ventilated")
$ws.Range("B51").Value = 22
$ws.Range("AB51").Value = 228
$ws = $wb.Worksheets.Item("Released")
$ws.Range("B51").Value = 400
$ws.Range("AB51").Value = 3540
